$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill data rows 2-38: I = 1 (constant), J = same value as H (IP)
for ($r = 2; $r -le 38; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
